# Brix/Eigenfertigung "Eigen_Glasgel" (glass-railing) sheet:
#   - The internal helper row "Intern: Rechn-Laenge" / "L_Calc" / "max(L, 1.0)"
#     is removed.
#   - Every formula that referenced the now-removed L_Calc variable is
#     rewritten to use the inline expression max(L, 1.0) instead.
#   - The length (L) and corner-count (Ecken) number inputs get a 2-decimal
#     number format, matching the other numeric input cell on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eigen_Glasgel")

# Remove the "Intern: Rechn-Länge / L_Calc / max(L, 1.0)" row entirely -
# this shifts every row below it (old rows 10-12) up by one.
$ws.Rows.Item(9).Delete()

# Rewrite the two formulas that used to reference L_Calc so they use the
# inline max(L, 1.0) expression instead (L_Calc no longer exists).
$ws.Range("E9").Value = "math.ceil(max(L, 1.0) / 1.3)"
$ws.Range("E11").Value = "(max(L, 1.0) * 0.85 * P_Glas) + (N_Steher * (P_Steher + F_Montage)) + (max(L, 1.0) * P_Handlauf) + ((N_Felder * 4 * P_Klem) + (Ecken * 4 * 5.0)) + (max(L, 1.0) * 65)"

# Give the length (D2) and corner-count (D3) inputs the same 2-decimal
# number format already used by the klemmhalter-count input (D8).
$ws.Range("D2").NumberFormat = "0.00"
$ws.Range("D3").NumberFormat = "0.00"

# Move the active selection to the last data cell, matching the saved state.
$ws.Activate() | Out-Null
$ws.Range("E11").Select() | Out-Null
